# Update cryptos list values (prices and 1h volume %) to reflect latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.996.52'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '3.386.07'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'571.57"
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = "'142.08"
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.474"
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').Value = '3.967.25'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').Value = "'27.87"
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '3.389.48'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '61.114.26'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  -3.09%  '
$ws.Range('D19').Value = "'13.59"
$ws.Range('E19').Value = '  -3.76%  '
$ws.Range('D20').Value = "'8.91"
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').Value = "'383.79"
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').Value = "'75.18"
$ws.Range('E22').Value = '  +2.76%  '
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').Value = '3.526.14'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = "'7.27"
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'7.97"
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.15"
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('D34').Value = "'23.22"
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').Value = "'6.95"
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('D36').Value = "'166.20"
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '3.420.01'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('D40').Value = "'0.0767"
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('D41').Value = "'26.85"
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = "'0.779"
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = "'4.37"
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('D46').Value = "'1.12"
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '2.453.08'
$ws.Range('E47').Value = '  -2.74%  '
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = "'6.71"
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = "'2.15"
$ws.Range('E50').Value = '  +10.85%  '
$ws.Range('D51').Value = "'0.0263"
$ws.Range('E51').Value = '  +1.50%  '
